$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This edit implements "Generate Report for handoff": the source file that
# used to be be339719-230e-4dd3-af51-baa301c6250b.md was regenerated with a
# new id (2f60501b-4111-4a8d-8abc-12f3957ff607.md) and a fresh handoff
# (new content hash + new handoff timestamps), and the
# 20a22625-7e2a-45b8-acaf-362483d04280.md row (status "Handoff transform
# failed") dropped out of the report entirely - the .localization-config row
# moves up to take its place.
# ---------------------------------------------------------------------------

$oldId = "be339719-230e-4dd3-af51-baa301c6250b"
$newId = "2f60501b-4111-4a8d-8abc-12f3957ff607"

$oldZhXlf = "$oldId.3a187063afe349459b7c35fdd3d89dbe3b2eac74.zh-cn.xlf"
$newZhXlf = "$newId.e63f8695803f223742439dd7b3a570b95cea9973.zh-cn.xlf"

$oldDeXlf = "$oldId.3a187063afe349459b7c35fdd3d89dbe3b2eac74.de-de.xlf"
$newDeXlf = "$newId.e63f8695803f223742439dd7b3a570b95cea9973.de-de.xlf"

$newZhTime = "2016-01-17 16:32:51"
$newDeTime = "2016-01-17 16:33:02"

function Fix-RowRemoval($ws, $lastCols) {
    # Row 3 (20a22625-....md / "Handoff transform failed") is being removed
    # from the report; row 4 (.localization-config) shifts up to row 3.
    foreach ($col in $lastCols) {
        $src = $ws.Range("$col 4".Replace(" ", ""))
        $dst = $ws.Range("$col 3".Replace(" ", ""))
        $dst.Value2 = $src.Value2
    }

    $lastCol = $lastCols[$lastCols.Length - 1]
    $ws.Range("A4:$lastCol 4".Replace(" ", "")).Clear()

    # The hyperlink that used to sit on A3 (20a22625...md) now has to show
    # the .localization-config row that was copied into row 3.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$A$3') {
            $h.TextToDisplay = ".localization-config"
        }
    }

    # The hyperlink that used to sit on A4 is now orphaned (row 4 is empty).
    $orphans = @()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$A$4') {
            $orphans += $h
        }
    }
    foreach ($h in $orphans) { $h.Delete() }
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "$newId.md"
foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "$newId.md"
    }
}

Fix-RowRemoval $wsOverview @("A", "B", "C")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value2 = "$newId.md"
$wsZh.Range("C2").Value2 = $newZhXlf
$wsZh.Range("D2").Value2 = $newZhTime

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newId.md"
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = $newZhXlf
    }
}

Fix-RowRemoval $wsZh @("A", "B", "C", "D", "E", "F", "G", "H", "I")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value2 = "$newId.md"
$wsDe.Range("C2").Value2 = $newDeXlf
$wsDe.Range("D2").Value2 = $newDeTime

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newId.md"
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = $newDeXlf
    }
}

Fix-RowRemoval $wsDe @("A", "B", "C", "D", "E", "F", "G", "H", "I")
